# "se sumaron mas imagenes al catalogo" — add more product codes to the
# MEJORAR catalog sheet.
#
# 1) "evol5530" was a duplicate of "EVOL5530" (case only) sitting at A37 —
#    remove that row, shifting everything below it up by one.
# 2) Append 9 new product codes at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEJORAR")

# Remove the stray duplicate "evol5530" row.
$ws.Rows.Item(37).Delete()

# Find the first empty row at the bottom of column A and append the new codes.
$newCodes = @(
    "EVOL0144",
    "EVOL5100",
    "EVOL0340",
    "EVOL0108",
    "EVOL1200",
    "EVOL1631",
    "EVOL2135",
    "EVOL0135",
    "EVOL1208"
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
foreach ($code in $newCodes) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $code
}

# Update the active selection to sit just past the new last row, matching
# where Excel would leave the cursor after appending these rows.
$ws.Range("A" + ($lastRow + 1)).Select()
